$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, pushing existing rows 141-199 down to 142-200.
# Excel automatically copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with its data.
$ws.Range("A141").Value = 11
$ws.Range("B141").Value = "Vega Monumental Concepción"
$ws.Range("C141").Value = "Bíobío"
$ws.Range("D141").Value = 44992
$ws.Range("E141").Value = 8
$ws.Range("F141").Value = 100112032
$ws.Range("G141").Value = "Zapallo italiano"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 100
$ws.Range("K141").Value = 8000
$ws.Range("L141").Value = 8500
$ws.Range("M141").Value = 8250
$ws.Range("N141").Value = "$/caja 50 unidades"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 165
$ws.Range("Q141").Value = 50
$ws.Range("R141").Value = "Hortaliza"
